$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5274.75
$ws.Range("I76").Value = 3699.6667
$ws.Range("K76").Value = 3699.6667
$ws.Range("M76").Value = -3384.6667
$ws.Range("H79").Value = 5274.75
$ws.Range("I79").Value = 3699.6667
$ws.Range("K79").Value = 3699.6667
$ws.Range("M79").Value = -2607.6667
$ws.Range("H80").Value = 1551.4736
$ws.Range("I80").Value = 332.75
$ws.Range("K80").Value = 998.25
$ws.Range("M80").Value = -0.25
$ws.Range("H83").Value = 1551.4736
$ws.Range("I83").Value = 332.75
$ws.Range("K83").Value = 2994.75
$ws.Range("M83").Value = 1997.25
$ws.Range("H92").Value = 1348.5
$ws.Range("I92").Value = 1741.5
$ws.Range("K92").Value = 1741.5
$ws.Range("M92").Value = -493.5
$ws.Range("H98").Value = 951.55554
$ws.Range("I98").Value = 858.125
$ws.Range("J98").Value = 1699
$ws.Range("K98").Value = 858.125
$ws.Range("L98").Value = 1699
$ws.Range("M98").Value = 639.875
$ws.Range("N98").Value = -4695
$ws.Range("H107").Value = 382
$ws.Range("I107").Value = 382
$ws.Range("K107").Value = 382
$ws.Range("M107").Value = 1538
$ws.Range("H122").Value = 951.55554
$ws.Range("I122").Value = 858.125
$ws.Range("J122").Value = 1699
$ws.Range("K122").Value = 2574.375
$ws.Range("L122").Value = 5097
$ws.Range("M122").Value = -124.375
$ws.Range("N122").Value = -9997
$ws.Range("H123").Value = 24999
$ws.Range("J123").Value = 24999
$ws.Range("L123").Value = 24999
$ws.Range("N123").Value = -34799
$ws.Range("H132").Value = 2189.5
$ws.Range("I132").Value = 2189.5
$ws.Range("K132").Value = 6568.5
$ws.Range("M132").Value = -4038.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1165.3334
$ws.Range("I102").Value = 1165.3334
$ws.Range("K102").Value = 1165.3334
$ws.Range("M102").Value = 456.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 585.3333
$ws.Range("I22").Value = 165
$ws.Range("K22").Value = 165
$ws.Range("M22").Value = 8
$ws.Range("H32").Value = 185555
$ws.Range("J32").Value = 185555
$ws.Range("L32").Value = 185555
$ws.Range("N32").Value = -186323
$ws.Range("H33").Value = 20000
$ws.Range("I33").Value = 15000
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -14664
$ws.Range("N33").Value = -25672
$ws.Range("H36").Value = 2997
$ws.Range("I36").Value = 1004.5
$ws.Range("J36").Value = 4989.5
$ws.Range("K36").Value = 1004.5
$ws.Range("L36").Value = 4989.5
$ws.Range("M36").Value = -470.5
$ws.Range("N36").Value = -6057.5
$ws.Range("H134").Value = 2081.55
$ws.Range("I134").Value = 2019.4706
$ws.Range("K134").Value = 6058.4118
$ws.Range("M134").Value = -3523.4118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 939.4
$ws.Range("J5").Value = 530.5
$ws.Range("L5").Value = 530.5
$ws.Range("N5").Value = -754.5
$ws.Range("H7").Value = 94.8
$ws.Range("I7").Value = 99
$ws.Range("J7").Value = 93.75
$ws.Range("K7").Value = 99
$ws.Range("L7").Value = 93.75
$ws.Range("M7").Value = 14
$ws.Range("N7").Value = -319.75
$ws.Range("H12").Value = 7709.75
$ws.Range("I12").Value = 8136.6
$ws.Range("K12").Value = 8136.6
$ws.Range("M12").Value = -7966.6
$ws.Range("H15").Value = 948
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H62").Value = 2795.2
$ws.Range("J62").Value = 2921.75
$ws.Range("L62").Value = 2921.75
$ws.Range("N62").Value = -4169.75
$ws.Range("H65").Value = 2795.2
$ws.Range("J65").Value = 2921.75
$ws.Range("L65").Value = 14608.75
$ws.Range("N65").Value = -20848.75
$ws.Range("H141").Value = 153428.56
$ws.Range("J141").Value = 153428.56
$ws.Range("L141").Value = 153428.56
$ws.Range("N141").Value = -163788.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 147.66667
$ws.Range("I17").Value = 121.5
$ws.Range("K17").Value = 364.5
$ws.Range("M17").Value = -195.5
$ws.Range("H50").Value = 509
$ws.Range("I50").Value = 420
$ws.Range("K50").Value = 1260
$ws.Range("M50").Value = -779
$ws.Range("H53").Value = 509
$ws.Range("I53").Value = 420
$ws.Range("K53").Value = 1260
$ws.Range("M53").Value = -779
$ws.Range("H75").Value = 10000
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H78").Value = 10000
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H121").Value = 14711.818
$ws.Range("J121").Value = 5180
$ws.Range("L121").Value = 15540
$ws.Range("N121").Value = -18160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7100
$ws.Range("J10").Value = 6900
$ws.Range("L10").Value = 6900
$ws.Range("N10").Value = -7238
$ws.Range("H70").Value = 5066.375
$ws.Range("J70").Value = 5339
$ws.Range("L70").Value = 5339
$ws.Range("N70").Value = -5879
$ws.Range("H73").Value = 5066.375
$ws.Range("J73").Value = 5339
$ws.Range("L73").Value = 5339
$ws.Range("N73").Value = -7211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2529.7856
$ws.Range("J46").Value = 3481.8572
$ws.Range("L46").Value = 3481.8572
$ws.Range("N46").Value = -3857.8572
$ws.Range("H136").Value = 3710
$ws.Range("I136").Value = 2539.6
$ws.Range("K136").Value = 7618.799999999999
$ws.Range("M136").Value = -5068.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2135.6428
$ws.Range("I132").Value = 2216.5833
$ws.Range("J132").Value = 1650
$ws.Range("K132").Value = 6649.749899999999
$ws.Range("L132").Value = 4950
$ws.Range("M132").Value = -4119.749899999999
$ws.Range("N132").Value = -10010
$ws.Range("H136").Value = 3170.5625
$ws.Range("I136").Value = 3540.8462
$ws.Range("K136").Value = 10622.5386
$ws.Range("M136").Value = -8072.5386
